$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 111815516
$ws.Range("B2").Value = 89423
$ws.Range("E2").Value = 5432
$ws.Range("F2").Value = 'Granticka'
$ws.Range("G2").Value = 'Porodaedalea chrysoloma'
$ws.Range("H2").Value = '(Fr.) Fiasson & Niemelä'
$ws.Range("K2").ClearContents()
$ws.Range("L2").ClearContents()
$ws.Range("M2").ClearContents()
$ws.Range("N2").ClearContents()
$ws.Range("Q2").Value = 458289.5512131723
$ws.Range("R2").Value = 7054475.069158822
$ws.Range("AC2").ClearContents()
$ws.Range("A3").Value = 111815507
$ws.Range("Q3").Value = 458151.5539710881
$ws.Range("R3").Value = 7054482.225765129
$ws.Range("AC3").Value = 'ringhack gamla'
$ws.Range("A4").Value = 111815515
$ws.Range("B4").Value = 89423
$ws.Range("E4").Value = 5432
$ws.Range("F4").Value = 'Granticka'
$ws.Range("G4").Value = 'Porodaedalea chrysoloma'
$ws.Range("H4").Value = '(Fr.) Fiasson & Niemelä'
$ws.Range("K4").ClearContents()
$ws.Range("L4").ClearContents()
$ws.Range("M4").ClearContents()
$ws.Range("N4").ClearContents()
$ws.Range("Q4").Value = 458161.9437607233
$ws.Range("R4").Value = 7054459.400503729
$ws.Range("AC4").ClearContents()
$ws.Range("A6").Value = 111815519
$ws.Range("B6").Value = 77515
$ws.Range("E6").Value = 6425
$ws.Range("F6").Value = 'Garnlav'
$ws.Range("G6").Value = 'Alectoria sarmentosa'
$ws.Range("H6").Value = '(Ach.) Ach.'
$ws.Range("Q6").Value = 458215.7474518137
$ws.Range("R6").Value = 7054621.063481365
$ws.Range("A7").Value = 111815512
$ws.Range("B7").Value = 56398
$ws.Range("E7").Value = 100109
$ws.Range("F7").Value = 'Tretåig hackspett'
$ws.Range("G7").Value = 'Picoides tridactylus'
$ws.Range("H7").Value = '(Linnaeus, 1758)'
# K7: empty placeholder cell (no representable value) - left blank
# L7: empty placeholder cell (no representable value) - left blank
# M7: empty placeholder cell (no representable value) - left blank
# N7: empty placeholder cell (no representable value) - left blank
$ws.Range("Q7").Value = 458154.6107204149
$ws.Range("R7").Value = 7054646.336103803
$ws.Range("AC7").Value = 'ringhack'
$ws.Range("A8").Value = 111815514
$ws.Range("B8").Value = 89423
$ws.Range("E8").Value = 5432
$ws.Range("F8").Value = 'Granticka'
$ws.Range("G8").Value = 'Porodaedalea chrysoloma'
$ws.Range("H8").Value = '(Fr.) Fiasson & Niemelä'
$ws.Range("K8").ClearContents()
$ws.Range("L8").ClearContents()
$ws.Range("M8").ClearContents()
$ws.Range("N8").ClearContents()
$ws.Range("Q8").Value = 458153.7808649908
$ws.Range("R8").Value = 7054482.19637617
$ws.Range("AC8").ClearContents()
$ws.Range("A9").Value = 111815513
$ws.Range("B9").Value = 56398
$ws.Range("E9").Value = 100109
$ws.Range("F9").Value = 'Tretåig hackspett'
$ws.Range("G9").Value = 'Picoides tridactylus'
$ws.Range("H9").Value = '(Linnaeus, 1758)'
# K9: empty placeholder cell (no representable value) - left blank
# L9: empty placeholder cell (no representable value) - left blank
# M9: empty placeholder cell (no representable value) - left blank
# N9: empty placeholder cell (no representable value) - left blank
$ws.Range("Q9").Value = 458173.7327805056
$ws.Range("R9").Value = 7054711.474791372
$ws.Range("AC9").Value = 'ringhack gamla'
$ws.Range("A10").Value = 111815517
$ws.Range("Q10").Value = 458250.8216980004
$ws.Range("R10").Value = 7054375.482693202
$ws.Range("A11").Value = 111815508
$ws.Range("B11").Value = 56398
$ws.Range("E11").Value = 100109
$ws.Range("F11").Value = 'Tretåig hackspett'
$ws.Range("G11").Value = 'Picoides tridactylus'
$ws.Range("H11").Value = '(Linnaeus, 1758)'
# K11: empty placeholder cell (no representable value) - left blank
# L11: empty placeholder cell (no representable value) - left blank
# M11: empty placeholder cell (no representable value) - left blank
# N11: empty placeholder cell (no representable value) - left blank
$ws.Range("Q11").Value = 458162.4570845839
$ws.Range("R11").Value = 7054329.489790585
$ws.Range("AC11").Value = 'ringhack'
$ws.Range("A12").Value = 111815518
$ws.Range("Q12").Value = 458250.901553072
$ws.Range("R12").Value = 7054618.376188213
